$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "Hoàng      Việt Bách"
[void]$ws.Range("C8").Select()
